$wb = $excel.ActiveWorkbook

# Reference to the existing first sheet so the new sheet is inserted right after it.
$ws1 = $wb.Worksheets.Item(1)

# Add the new "warning_position" sheet right after Sheet1 (this also makes it the
# active sheet/tab, matching the added activeTab="1" on the workbookView).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "warning_position"

# Populate the new sheet's data (B3:B6).
$ws2.Range("B3").Value = "var1"
$ws2.Range("B4").Value = 1234
$ws2.Range("B5").Value = 2345
$ws2.Range("B6").Value = "wtf"

# Match the saved selection/active cell on the new sheet.
$ws2.Range("B6").Select() | Out-Null
